# Add a conditional block around the existing {questionTypeName} list item so
# the worksheet PDF template shows a simple "{questionTypeName}" line when
# there is only a single question type, and falls back to the original
# bulleted "{questionTypeName}" list item otherwise.
#
#   {#isSingleQuestionType}
#   {questionTypeName}
#   {/isSingleQuestionType}
#   {^isSingleQuestionType}
#   <original bulleted "{questionTypeName}" paragraph, unchanged>
#   {/isSingleQuestionType}

$d = $word.ActiveDocument

$wdNamespace = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ParagraphAfter($para, [string]$xml) {
    # Creates a brand-new empty paragraph right after $para, then fills it
    # with the supplied WordprocessingML paragraph XML (pPr/runs/etc.).
    # Returns the Paragraph object for the newly-inserted paragraph.
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $insertionPoint = $newPara.Range.Duplicate()
    $insertionPoint.Collapse(1)
    $insertionPoint.InsertXML($xml)
    return $para.Next()
}

function Find-ParagraphByText([string]$text) {
    # Paragraph objects captured before an earlier InsertXML/InsertParagraphAfter
    # can go stale (they keep their old character offsets instead of floating
    # with the content), so anchors are always (re-)located fresh right before
    # they are used.
    $fr = $d.Content
    $fr.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $fr.Paragraphs(1)
}

# --- Hunk 1: four new paragraphs right after "{#questionTypes}" -----------------

$xmlOpenSingle = "<w:p $wdNamespace><w:pPr><w:pStyle w:val=`"Normal`"/><w:widowControl w:val=`"0`"/><w:ind w:left=`"720`" w:firstLine=`"0`"/></w:pPr><w:r><w:rPr/><w:t>{#isSingleQuestionType}</w:t></w:r></w:p>"

$xmlSingleName = "<w:p $wdNamespace><w:pPr><w:pStyle w:val=`"Normal`"/><w:widowControl w:val=`"0`"/><w:ind w:left=`"720`"/><w:rPr><w:rFonts w:ascii=`"Open Sans`" w:hAnsi=`"Open Sans`" w:eastAsia=`"Open Sans`" w:cs=`"Open Sans`"/><w:b w:val=`"1`"/><w:bCs w:val=`"1`"/><w:i w:val=`"0`"/><w:iCs w:val=`"0`"/><w:color w:val=`"19A27D`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Open Sans`" w:hAnsi=`"Open Sans`" w:eastAsia=`"Open Sans`" w:cs=`"Open Sans`"/><w:b w:val=`"1`"/><w:bCs w:val=`"1`"/><w:i w:val=`"0`"/><w:iCs w:val=`"0`"/><w:color w:val=`"19A27D`"/></w:rPr><w:t>{questionTypeName}</w:t></w:r></w:p>"

$xmlCloseSingle = "<w:p $wdNamespace><w:pPr><w:pStyle w:val=`"Normal`"/><w:widowControl w:val=`"0`"/><w:ind w:left=`"720`" w:firstLine=`"0`"/></w:pPr><w:r><w:rPr/><w:t>{/isSingleQuestionType}</w:t></w:r></w:p>"

$xmlOpenNotSingle = "<w:p $wdNamespace><w:pPr><w:pStyle w:val=`"Normal`"/><w:widowControl w:val=`"0`"/><w:ind w:left=`"720`" w:firstLine=`"0`"/></w:pPr><w:r><w:rPr/><w:t>{^isSingleQuestionType}</w:t></w:r></w:p>"

$cursor = Find-ParagraphByText("{#questionTypes}")
$cursor = Insert-ParagraphAfter $cursor $xmlOpenSingle
$cursor = Insert-ParagraphAfter $cursor $xmlSingleName
$cursor = Insert-ParagraphAfter $cursor $xmlCloseSingle
$cursor = Insert-ParagraphAfter $cursor $xmlOpenNotSingle

# --- Hunk 2: one new paragraph right after the original bulleted ----------------
# --- "{questionTypeName}" list paragraph, closing the "^isSingleQuestionType" ---
# --- block before "{instruction}" ------------------------------------------------

$instructionPara = Find-ParagraphByText("{instruction}")
$questionTypeNameListPara = $instructionPara.Previous()

$xmlCloseNotSingle = "<w:p $wdNamespace><w:pPr><w:pStyle w:val=`"Normal`"/><w:widowControl w:val=`"0`"/><w:ind w:left=`"720`" w:firstLine=`"0`"/></w:pPr><w:r><w:rPr/><w:t>{/isSingleQuestionType}</w:t></w:r></w:p>"

Insert-ParagraphAfter $questionTypeNameListPara $xmlCloseNotSingle | Out-Null

Write-Host "Inserted isSingleQuestionType condition around questionTypeName block."
